$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1030.238
$ws.Range("I98").Value = 970.1579
$ws.Range("J98").Value = 1601
$ws.Range("K98").Value = 970.1579
$ws.Range("L98").Value = 1601
$ws.Range("M98").Value = 527.8421
$ws.Range("N98").Value = -4597
$ws.Range("H116").Value = 8312.958000000001
$ws.Range("I116").Value = 7815.6665
$ws.Range("J116").Value = 8478.723
$ws.Range("K116").Value = 7815.6665
$ws.Range("L116").Value = 8478.723
$ws.Range("M116").Value = -4373.6665
$ws.Range("N116").Value = -15362.723
$ws.Range("H122").Value = 1030.238
$ws.Range("I122").Value = 970.1579
$ws.Range("J122").Value = 1601
$ws.Range("K122").Value = 2910.4737
$ws.Range("L122").Value = 4803
$ws.Range("M122").Value = -460.4737
$ws.Range("N122").Value = -9703
$ws.Range("H138").Value = 3445.9404
$ws.Range("I138").Value = 1247.6957
$ws.Range("J138").Value = 4274.787
$ws.Range("K138").Value = 3743.0871
$ws.Range("L138").Value = 12824.361
$ws.Range("M138").Value = 1396.9129
$ws.Range("N138").Value = -23104.361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 14509.85
$ws.Range("I45").Value = 10399.111
$ws.Range("K45").Value = 10399.111
$ws.Range("M45").Value = -10022.111
$ws.Range("H74").Value = 23819466
$ws.Range("I74").Value = 7066.8335
$ws.Range("J74").Value = 55569332
$ws.Range("K74").Value = 7066.8335
$ws.Range("L74").Value = 55569332
$ws.Range("M74").Value = -6192.8335
$ws.Range("N74").Value = -55571080
$ws.Range("H77").Value = 23819466
$ws.Range("I77").Value = 7066.8335
$ws.Range("J77").Value = 55569332
$ws.Range("K77").Value = 35334.1675
$ws.Range("L77").Value = 277846660
$ws.Range("M77").Value = -30966.1675
$ws.Range("N77").Value = -277855396
$ws.Range("H122").Value = 1010.3077
$ws.Range("I122").Value = 927.8333
$ws.Range("K122").Value = 2783.4999
$ws.Range("M122").Value = -333.4998999999998
$ws.Range("H132").Value = 2145.7144
$ws.Range("I132").Value = 1608.4722
$ws.Range("J132").Value = 5369.1665
$ws.Range("K132").Value = 4825.4166
$ws.Range("L132").Value = 16107.4995
$ws.Range("M132").Value = -2295.4166
$ws.Range("N132").Value = -21167.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 434974.28
$ws.Range("I22").Value = 596.8182
$ws.Range("J22").Value = 532487.5600000001
$ws.Range("K22").Value = 596.8182
$ws.Range("L22").Value = 532487.5600000001
$ws.Range("M22").Value = -423.8182
$ws.Range("N22").Value = -532833.5600000001
$ws.Range("H86").Value = 5684.5835
$ws.Range("I86").Value = 3277.75
$ws.Range("J86").Value = 10498.25
$ws.Range("K86").Value = 3277.75
$ws.Range("L86").Value = 10498.25
$ws.Range("M86").Value = -2154.75
$ws.Range("N86").Value = -12744.25
$ws.Range("H89").Value = 5684.5835
$ws.Range("I89").Value = 3277.75
$ws.Range("J89").Value = 10498.25
$ws.Range("K89").Value = 16388.75
$ws.Range("L89").Value = 52491.25
$ws.Range("M89").Value = -10772.75
$ws.Range("N89").Value = -63723.25
$ws.Range("H94").Value = 1280.1578
$ws.Range("I94").Value = 1197.8667
$ws.Range("K94").Value = 1197.8667
$ws.Range("M94").Value = -746.8667
$ws.Range("H99").Value = 5879.9
$ws.Range("I99").Value = 2258.1667
$ws.Range("J99").Value = 11312.5
$ws.Range("K99").Value = 2258.1667
$ws.Range("L99").Value = 11312.5
$ws.Range("M99").Value = -760.1667000000002
$ws.Range("N99").Value = -14308.5
$ws.Range("H107").Value = 2800
$ws.Range("I107").Value = 3700
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 3700
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = -1780
$ws.Range("N107").Value = -5740
$ws.Range("H134").Value = 1850.5807
$ws.Range("I134").Value = 1663.0741
$ws.Range("K134").Value = 4989.2223
$ws.Range("M134").Value = -2454.2223
$ws.Range("H138").Value = 108842
$ws.Range("J138").Value = 108842
$ws.Range("L138").Value = 108842
$ws.Range("N138").Value = -119122

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 21915
$ws.Range("I33").Value = 1166.6666
$ws.Range("K33").Value = 1166.6666
$ws.Range("M33").Value = -787.6666
$ws.Range("H134").Value = 13948.432
$ws.Range("I134").Value = 14603.269
$ws.Range("K134").Value = 43809.807
$ws.Range("M134").Value = -41274.807

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1924315.6
$ws.Range("J68").Value = 2501302.5
$ws.Range("L68").Value = 7503907.5
$ws.Range("N68").Value = -7505529.5
$ws.Range("H71").Value = 1924315.6
$ws.Range("J71").Value = 2501302.5
$ws.Range("L71").Value = 22511722.5
$ws.Range("N71").Value = -22519834.5
$ws.Range("H86").Value = 749.25
$ws.Range("J86").Value = 750
$ws.Range("L86").Value = 2250
$ws.Range("N86").Value = -4622
$ws.Range("H89").Value = 749.25
$ws.Range("J89").Value = 750
$ws.Range("L89").Value = 6750
$ws.Range("N89").Value = -18606

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 372.68
$ws.Range("J107").Value = 500.44446
$ws.Range("L107").Value = 500.44446
$ws.Range("N107").Value = -4340.44446
$ws.Range("H113").Value = 3727.8572
$ws.Range("I113").Value = 3092.4
$ws.Range("K113").Value = 3092.4
$ws.Range("M113").Value = -922.4000000000001
$ws.Range("H126").Value = 22782.53
$ws.Range("I126").Value = 28031
$ws.Range("K126").Value = 84093
$ws.Range("M126").Value = -81623
$ws.Range("H132").Value = 27851.762
$ws.Range("I132").Value = 28994.35
$ws.Range("K132").Value = 86983.04999999999
$ws.Range("M132").Value = -84453.04999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8667.333000000001
$ws.Range("I7").Value = 10667.5
$ws.Range("K7").Value = 10667.5
$ws.Range("M7").Value = -10555.5
$ws.Range("H55").Value = 1231.8182
$ws.Range("H61").Value = 5702.2856
$ws.Range("I61").Value = 5114.2
$ws.Range("J61").Value = 7172.5
$ws.Range("K61").Value = 5114.2
$ws.Range("L61").Value = 7172.5
$ws.Range("M61").Value = -4912.2
$ws.Range("N61").Value = -7576.5
$ws.Range("H82").Value = 4697.55
$ws.Range("J82").Value = 5539.3
$ws.Range("L82").Value = 5539.3
$ws.Range("N82").Value = -6261.3
$ws.Range("H85").Value = 4697.55
$ws.Range("J85").Value = 5539.3
$ws.Range("L85").Value = 5539.3
$ws.Range("N85").Value = -8035.3
$ws.Range("H113").Value = 5702.2856
$ws.Range("I113").Value = 5114.2
$ws.Range("J113").Value = 7172.5
$ws.Range("K113").Value = 5114.2
$ws.Range("L113").Value = 7172.5
$ws.Range("M113").Value = -2944.2
$ws.Range("N113").Value = -11512.5
$ws.Range("H122").Value = 4481.6
$ws.Range("I122").Value = 3825.4443
$ws.Range("J122").Value = 5465.8335
$ws.Range("K122").Value = 11476.3329
$ws.Range("L122").Value = 16397.5005
$ws.Range("M122").Value = -9026.332900000001
$ws.Range("N122").Value = -21297.5005
$ws.Range("H126").Value = 8667.333000000001
$ws.Range("I126").Value = 10667.5
$ws.Range("K126").Value = 32002.5
$ws.Range("M126").Value = -29532.5
$ws.Range("H136").Value = 4158.6665
$ws.Range("I136").Value = 3635.5334
$ws.Range("J136").Value = 5466.5
$ws.Range("K136").Value = 10906.6002
$ws.Range("L136").Value = 16399.5
$ws.Range("M136").Value = -8356.600199999999
$ws.Range("N136").Value = -21499.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8150.8423
$ws.Range("I81").Value = 12059.111
$ws.Range("J81").Value = 4633.4
$ws.Range("K81").Value = 24118.222
$ws.Range("L81").Value = 9266.799999999999
$ws.Range("M81").Value = -23057.222
$ws.Range("N81").Value = -11388.8
$ws.Range("H84").Value = 8150.8423
$ws.Range("I84").Value = 12059.111
$ws.Range("J84").Value = 4633.4
$ws.Range("K84").Value = 120591.11
$ws.Range("L84").Value = 46334
$ws.Range("M84").Value = -115287.11
$ws.Range("N84").Value = -56942
$ws.Range("H122").Value = 2527.2856
$ws.Range("I122").Value = 2410.56
$ws.Range("K122").Value = 7231.68
$ws.Range("M122").Value = -4781.68
$ws.Range("H126").Value = 3033.08
$ws.Range("I126").Value = 2546.1
$ws.Range("J126").Value = 4981
$ws.Range("K126").Value = 7638.299999999999
$ws.Range("L126").Value = 14943
$ws.Range("M126").Value = -5168.299999999999
$ws.Range("N126").Value = -19883
$ws.Range("H132").Value = 3083.6135
$ws.Range("I132").Value = 3283.6843
$ws.Range("K132").Value = 9851.052899999999
$ws.Range("M132").Value = -7321.052899999999
